# The commit adds a "status" column (with a "pass" value for each data row)
# to the RegdTestData sheet, mirroring the existing "status" column F.
# New cells: G1/H1 = "status" header, G2:G7 = "pass" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegdTestData")

# Header cells - same text & style treatment as the existing "status" header in F1
$ws.Range("G1").Value = "status"
$ws.Range("H1").Value = "status"
$ws.Range("G1").Style = $ws.Range("F1").Style
$ws.Range("H1").Style = $ws.Range("F1").Style

# Data cells - same text as the existing "pass" values in column F
$ws.Range("G2").Value = "pass"
$ws.Range("G3").Value = "pass"
$ws.Range("G4").Value = "pass"
$ws.Range("G5").Value = "pass"
$ws.Range("G6").Value = "pass"
$ws.Range("G7").Value = "pass"

# Match column width behaviour of the adjacent (same-content) column F
$ws.Columns("G:G").AutoFit()
